$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # "Pieces"
$ws2 = $wb.Worksheets.Item(2)  # "SRLs"

# The linked query tables ("Consulta desde VRT_TRAKING") were refreshed
# against the external SQL Server source; with no data available the
# returned row sets are now empty, so the previously-fetched rows are
# cleared from both worksheets.
$ws1.Range("A6:M7").ClearContents()
$ws2.Range("A6:H7").ClearContents()

# The invoice start/end date cells (populated by the same refresh) are
# cleared as well.
$ws1.Range("K1:K2").ClearContents()
$ws2.Range("J1:J2").ClearContents()

# Selection / active-sheet bookkeeping left behind by the interactive
# session that triggered the refresh.
[void]$ws1.Range("L1").Select()
[void]$ws2.Activate()
[void]$ws2.Range("B15").Select()
